$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Invoice Number INV-333 "
$ws.Range("B2").Value = "From: DEMO Sliced Invoices Suite 5A-1204 123 Somewhere Street Your Citv 4Z 12345 "
$ws.Range("C2").Value = "To: Test Business 123 Somewhere St Melbourne, VIC 3000 testtest com "
$ws.Range("E2").Value = "test@test, admin@slicedinvoices, "
